$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with latest scraped values.
# Cells whose new text looks like a plain number (e.g. "243.52") must have their
# number format forced to Text ("@") first, otherwise Excel auto-converts the string
# into a numeric value instead of keeping it as the literal text used on the source site.

$ws.Range('D2').Value = '29.413.99'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.868.50'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.52'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7044'
$ws.Range('E6').Value = '  -2.28%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07962'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3137'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.47'
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07816'
$ws.Range('E11').Value = '  -4.74%  '
$ws.Range('D12').Value = '1.905.73'
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '93.80'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.172'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7031'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.489'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008615'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').Value = '29.534.61'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.28'
$ws.Range('E19').Value = '  +3.60%  '
$ws.Range('D20').Value = '2.156.76'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.658'
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1550'
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.006'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.46'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.508'
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('E30').Value = '  -2.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.259'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05274'
$ws.Range('E33').Value = '  -1.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.906'
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7585'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.184'
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01880'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').Value = '1.279.83'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.770'
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8962'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.89'
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.032'
$ws.Range('E43').Value = '  -6.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '70.95'
$ws.Range('E44').Value = '  -4.31%  '
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').Value = '2.047.08'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  -3.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.809'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.607'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.5181'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4298'
$ws.Range('E51').Value = '  -1.11%  '

Write-Host "Applied 89 cell updates"
